$d = $word.ActiveDocument

# Locate the paragraph that ends the presentation and contains the
# sentence we need to split into several runs.
$needle = "Si vous voulez en savoir plus je vous invite"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith($needle)) {
        $target = $p
    }
}

$paraStart = $target.Range.Start
$fullText = $target.Range.Text
$marker = "Télescope"
$splitOffset = $fullText.IndexOf($marker) + $marker.Length
$splitPos = $paraStart + $splitOffset

$insert1 = " qui ciblera plus précisément "
$insert2 = "Les défis informatiques et de traitement de données"

# --- Insert " qui ciblera plus précisément " right after "Télescope" ---
$gap1 = $d.Range($splitPos, $splitPos)
$gap1.InsertAfter($insert1)
$seg1Start = $splitPos
$seg1End = $seg1Start + $insert1.Length
# Force this new text into its own run instead of letting it merge back
# into the identically-formatted text around it.
$seg1 = $d.Range($seg1Start, $seg1End)
$seg1.Bold = 1
$seg1.Bold = 0

# --- Insert "Les défis informatiques et de traitement de données" next ---
$gap2 = $d.Range($seg1End, $seg1End)
$gap2.InsertAfter($insert2)
$seg2Start = $seg1End
$seg2End = $seg2Start + $insert2.Length
$seg2 = $d.Range($seg2Start, $seg2End)
$seg2.Bold = 1
$seg2.Bold = 0

# --- Keep the trailing ". Merci de votre écoute" as its own run too ---
$tailStart = $seg2End
$tailEnd = $target.Range.End
$tail = $d.Range($tailStart, $tailEnd)
$tail.Bold = 1
$tail.Bold = 0
